# Insert a new weekly record at the top of the Coliflor price table
# (row 444), pushing the existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 444:469 down to 445:470, carrying formatting with them.
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with the new weekly data point.
$ws.Range("A444").Value = 4
$ws.Range("B444").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C444").Value = "Los Lagos"
$ws.Range("D444").Value = 44939
$ws.Range("E444").Value = 10
$ws.Range("F444").Value = 100112008
$ws.Range("G444").Value = "Coliflor"
$ws.Range("H444").Value = "Sin especificar"
$ws.Range("I444").Value = "Primera"
$ws.Range("J444").Value = 1200
$ws.Range("K444").Value = 1600
$ws.Range("L444").Value = 1600
$ws.Range("M444").Value = 1600
$ws.Range("N444").Value = "$/unidad"
$ws.Range("O444").Value = "Región Metropolitana"
$ws.Range("P444").Value = 1600
$ws.Range("Q444").Value = 1
$ws.Range("R444").Value = "Hortaliza"
